$d = $word.ActiveDocument
$parts = $d.CustomXMLParts
Write-Output ("Count=" + $parts.Count)
try {
    $all = $parts.Item(1)
    Write-Output "got item 1"
} catch {
    Write-Output ("ERR: " + $_.Exception.Message)
}
try {
    $newPart = $parts.Add("<root/>")
    Write-Output "added new part"
    Write-Output ("Count now=" + $parts.Count)
} catch {
    Write-Output ("ERR add: " + $_.Exception.Message)
}
